$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.241005420684814
$ws.Range("B1").Value = 1.418902039527893
$ws.Range("C1").Value = 1.770377993583679
$ws.Range("D1").Value = 3.486965894699097
$ws.Range("E1").Value = 15
